# Update the "Corr/total marks" figures on the marksheet (quiz worksheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking" -> Right column (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Row 12 "Total" -> Right column (B12): 33 -> 55
$ws.Range("B12").Value = 55

# Row 12 "Total" -> Max column (E12): "17/84" -> "55/140"
$ws.Range("E12").Value = "55/140"
